$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 999.7143
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H43").Value = 4926.6665
$ws.Range("I43").Value = 2980
$ws.Range("J43").Value = 5900
$ws.Range("K43").Value = 2980
$ws.Range("L43").Value = 5900
$ws.Range("M43").Value = -2911
$ws.Range("N43").Value = -6038
$ws.Range("H129").Value = 815.3200000000001
$ws.Range("J129").Value = 867.75824
$ws.Range("L129").Value = 2603.27472
$ws.Range("N129").Value = -12603.27472
$ws.Range("H141").Value = 41103.77
$ws.Range("I141").Value = 57546.832
$ws.Range("J141").Value = 4106.875
$ws.Range("K141").Value = 172640.496
$ws.Range("L141").Value = 12320.625
$ws.Range("M141").Value = -167460.496
$ws.Range("N141").Value = -22680.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23336.166
$ws.Range("I37").Value = 9825
$ws.Range("J37").Value = 30091.75
$ws.Range("K37").Value = 9825
$ws.Range("L37").Value = 30091.75
$ws.Range("M37").Value = -9552
$ws.Range("N37").Value = -30637.75
$ws.Range("H55").Value = 25757
$ws.Range("J55").Value = 25757
$ws.Range("L55").Value = 25757
$ws.Range("N55").Value = -26387
$ws.Range("H80").Value = 42651.2
$ws.Range("J80").Value = 42651.2
$ws.Range("L80").Value = 42651.2
$ws.Range("N80").Value = -44647.2
$ws.Range("H83").Value = 42651.2
$ws.Range("J83").Value = 42651.2
$ws.Range("L83").Value = 127953.6
$ws.Range("N83").Value = -137937.6
$ws.Range("H103").Value = 34978.69
$ws.Range("J103").Value = 34978.69
$ws.Range("L103").Value = 34978.69
$ws.Range("N103").Value = -37322.69
$ws.Range("H110").Value = 1000.9167
$ws.Range("I110").Value = 1091.8334
$ws.Range("J110").Value = 910
$ws.Range("K110").Value = 1091.8334
$ws.Range("L110").Value = 910
$ws.Range("M110").Value = 953.1666
$ws.Range("N110").Value = -5000
$ws.Range("H132").Value = 2119.3157
$ws.Range("I132").Value = 812.4286
$ws.Range("J132").Value = 5778.6
$ws.Range("K132").Value = 2437.2858
$ws.Range("L132").Value = 17335.8
$ws.Range("M132").Value = 92.71420000000035
$ws.Range("N132").Value = -22395.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2325.625
$ws.Range("I86").Value = 1602.6666
$ws.Range("J86").Value = 2759.4
$ws.Range("K86").Value = 1602.6666
$ws.Range("L86").Value = 2759.4
$ws.Range("M86").Value = -479.6666
$ws.Range("N86").Value = -5005.4
$ws.Range("H89").Value = 2325.625
$ws.Range("I89").Value = 1602.6666
$ws.Range("J89").Value = 2759.4
$ws.Range("K89").Value = 8013.333000000001
$ws.Range("L89").Value = 13797
$ws.Range("M89").Value = -2397.333000000001
$ws.Range("N89").Value = -25029
$ws.Range("H128").Value = 1300
$ws.Range("I128").Value = 1300
$ws.Range("K128").Value = 3900
$ws.Range("M128").Value = -1410
$ws.Range("H134").Value = 3010.6365
$ws.Range("I134").Value = 1789.2222
$ws.Range("J134").Value = 8507
$ws.Range("K134").Value = 5367.6666
$ws.Range("L134").Value = 25521
$ws.Range("M134").Value = -2832.6666
$ws.Range("N134").Value = -30591

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3878.111
$ws.Range("I31").Value = 1788.8
$ws.Range("J31").Value = 5107.1177
$ws.Range("K31").Value = 1788.8
$ws.Range("L31").Value = 5107.1177
$ws.Range("M31").Value = -1493.8
$ws.Range("N31").Value = -5697.1177
$ws.Range("H34").Value = 3878.111
$ws.Range("I34").Value = 1788.8
$ws.Range("J34").Value = 5107.1177
$ws.Range("K34").Value = 1788.8
$ws.Range("L34").Value = 5107.1177
$ws.Range("M34").Value = -1586.8
$ws.Range("N34").Value = -5511.1177
$ws.Range("H99").Value = 6899613.5
$ws.Range("I99").Value = 18183396
$ws.Range("J99").Value = 3968.5
$ws.Range("K99").Value = 18183396
$ws.Range("L99").Value = 3968.5
$ws.Range("M99").Value = -18181898
$ws.Range("N99").Value = -6964.5
$ws.Range("H122").Value = 2778.3
$ws.Range("J122").Value = 5933.3335
$ws.Range("L122").Value = 17800.0005
$ws.Range("N122").Value = -22700.0005
$ws.Range("H126").Value = 6899613.5
$ws.Range("I126").Value = 18183396
$ws.Range("J126").Value = 3968.5
$ws.Range("K126").Value = 54550188
$ws.Range("L126").Value = 11905.5
$ws.Range("M126").Value = -54547718
$ws.Range("N126").Value = -16845.5
$ws.Range("H139").Value = 39000
$ws.Range("J139").Value = 39000
$ws.Range("L139").Value = 39000
$ws.Range("M139").Value = -49280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2088.9246
$ws.Range("I121").Value = 352.66666
$ws.Range("J121").Value = 2193.1
$ws.Range("K121").Value = 1057.99998
$ws.Range("L121").Value = 6579.299999999999
$ws.Range("M121").Value = 252.0000199999999
$ws.Range("N121").Value = -9199.299999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 21899
$ws.Range("J48").Value = 21899
$ws.Range("L48").Value = 21899
$ws.Range("N48").Value = -22869
$ws.Range("H70").Value = 5330.017
$ws.Range("I70").Value = 5056.073
$ws.Range("J70").Value = 5954
$ws.Range("K70").Value = 5056.073
$ws.Range("L70").Value = 5954
$ws.Range("M70").Value = -4786.073
$ws.Range("N70").Value = -6494
$ws.Range("H73").Value = 5330.017
$ws.Range("I73").Value = 5056.073
$ws.Range("J73").Value = 5954
$ws.Range("K73").Value = 5056.073
$ws.Range("L73").Value = 5954
$ws.Range("M73").Value = -4120.073
$ws.Range("N73").Value = -7826
$ws.Range("H102").Value = 2079.2942
$ws.Range("I102").Value = 1507
$ws.Range("J102").Value = 4750
$ws.Range("K102").Value = 1507
$ws.Range("L102").Value = 4750
$ws.Range("M102").Value = 115
$ws.Range("N102").Value = -7994
$ws.Range("H122").Value = 5400.9287
$ws.Range("I122").Value = 1701.5555
$ws.Range("J122").Value = 12059.8
$ws.Range("K122").Value = 5104.666499999999
$ws.Range("L122").Value = 36179.39999999999
$ws.Range("M122").Value = -2654.666499999999
$ws.Range("N122").Value = -41079.39999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7116.364
$ws.Range("I7").Value = 6885
$ws.Range("K7").Value = 6885
$ws.Range("M7").Value = -6773
$ws.Range("H18").Value = 19893.5
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 19893.5
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 19893.5
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -20237.5
$ws.Range("H20").Value = 6210.5293
$ws.Range("I20").Value = 2980
$ws.Range("J20").Value = 9844.875
$ws.Range("K20").Value = 2980
$ws.Range("L20").Value = 9844.875
$ws.Range("M20").Value = -2754
$ws.Range("N20").Value = -10296.875
$ws.Range("H122").Value = 8396.666999999999
$ws.Range("I122").Value = 6126.6665
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 18379.9995
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -15929.9995
$ws.Range("N122").Value = -36900.001
$ws.Range("H126").Value = 7116.364
$ws.Range("I126").Value = 6885
$ws.Range("K126").Value = 20655
$ws.Range("M126").Value = -18185

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3505.7646
$ws.Range("I122").Value = 2978.4285
$ws.Range("J122").Value = 5966.6665
$ws.Range("K122").Value = 8935.2855
$ws.Range("L122").Value = 17899.9995
$ws.Range("M122").Value = -6485.2855
$ws.Range("N122").Value = -22799.9995
